$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 197.88235
$ws.Range("J5").Value = 197.75
$ws.Range("L5").Value = 197.75
$ws.Range("N5").Value = -427.75

$ws.Range("H28").Value = 1258.2858
$ws.Range("I28").Value = 968.8333
$ws.Range("K28").Value = 968.8333
$ws.Range("M28").Value = -483.8333

$ws.Range("H32").Value = 11476
$ws.Range("J32").Value = 11845
$ws.Range("L32").Value = 11845
$ws.Range("N32").Value = -12497

$ws.Range("H49").Value = 124
$ws.Range("I49").Value = 199
$ws.Range("J49").Value = 49
$ws.Range("K49").Value = 597
$ws.Range("L49").Value = 147
$ws.Range("M49").Value = -461
$ws.Range("N49").Value = -419

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = $null

$ws.Range("H62").Value = 5622.5
$ws.Range("I62").Value = 6997.3335
$ws.Range("J62").Value = 1498
$ws.Range("K62").Value = 6997.3335
$ws.Range("L62").Value = 1498
$ws.Range("M62").Value = -6373.3335
$ws.Range("N62").Value = -2746

$ws.Range("H65").Value = 5622.5
$ws.Range("I65").Value = 6997.3335
$ws.Range("J65").Value = 1498
$ws.Range("K65").Value = 34986.6675
$ws.Range("L65").Value = 7490
$ws.Range("M65").Value = -31866.6675
$ws.Range("N65").Value = -13730

$ws.Range("H96").Value = 751.7778
$ws.Range("I96").Value = 709.4286
$ws.Range("K96").Value = 2128.2858
$ws.Range("M96").Value = -755.2857999999997

$ws.Range("H98").Value = 871.1667
$ws.Range("I98").Value = 306.75
$ws.Range("K98").Value = 306.75
$ws.Range("M98").Value = 1191.25

$ws.Range("H101").Value = 1580
$ws.Range("J101").Value = 2000
$ws.Range("L101").Value = 6000
$ws.Range("N101").Value = -9244

$ws.Range("H106").Value = 4999.6665
$ws.Range("I106").Value = 4999.6665
$ws.Range("K106").Value = 4999.6665
$ws.Range("M106").Value = -4368.6665

$ws.Range("H112").Value = 1600.1111
$ws.Range("J112").Value = 1600.1111
$ws.Range("L112").Value = 4800.3333
$ws.Range("N112").Value = -7016.3333

$ws.Range("H121").Value = 1852.1578
$ws.Range("J121").Value = 1949.7778
$ws.Range("L121").Value = 5849.3334
$ws.Range("N121").Value = -9343.3334

$ws.Range("H122").Value = 871.1667
$ws.Range("I122").Value = 306.75
$ws.Range("K122").Value = 920.25
$ws.Range("M122").Value = 1529.75

$ws.Range("H125").Value = 4239.75
$ws.Range("I125").Value = 3911.8572
$ws.Range("J125").Value = 4698.8
$ws.Range("K125").Value = 35206.7148
$ws.Range("L125").Value = 42289.2
$ws.Range("M125").Value = -32746.7148
$ws.Range("N125").Value = -47209.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 930.5
$ws.Range("I2").Value = 846.3077
$ws.Range("J2").Value = 1149.4
$ws.Range("K2").Value = 846.3077
$ws.Range("L2").Value = 1149.4
$ws.Range("M2").Value = -733.3077
$ws.Range("N2").Value = -1375.4

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").Value = $null

$ws.Range("H45").Value = 1429.4445
$ws.Range("I45").Value = 1429.4445
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1429.4445
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1052.4445
$ws.Range("N45").Value = $null

$ws.Range("H61").Value = 2598.2
$ws.Range("I61").Value = 2598.2
$ws.Range("K61").Value = 2598.2
$ws.Range("M61").Value = -2386.2

$ws.Range("H74").Value = 9998168
$ws.Range("I74").Value = 14281098
$ws.Range("K74").Value = 14281098
$ws.Range("M74").Value = -14280224

$ws.Range("H77").Value = 9998168
$ws.Range("I77").Value = 14281098
$ws.Range("K77").Value = 71405490
$ws.Range("M77").Value = -71401122

$ws.Range("H116").Value = 930.5
$ws.Range("I116").Value = 846.3077
$ws.Range("J116").Value = 1149.4
$ws.Range("K116").Value = 846.3077
$ws.Range("L116").Value = 1149.4
$ws.Range("M116").Value = 1447.6923
$ws.Range("N116").Value = -5737.4

$ws.Range("H136").Value = 2598.2
$ws.Range("I136").Value = 2598.2
$ws.Range("K136").Value = 7794.599999999999
$ws.Range("M136").Value = -5244.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 930.5
$ws.Range("I3").Value = 846.3077
$ws.Range("J3").Value = 1149.4
$ws.Range("K3").Value = 846.3077
$ws.Range("L3").Value = 1149.4
$ws.Range("M3").Value = -732.3077
$ws.Range("N3").Value = -1377.4

$ws.Range("H86").Value = 4392.1113
$ws.Range("I86").Value = 3622.2
$ws.Range("J86").Value = 5354.5
$ws.Range("K86").Value = 3622.2
$ws.Range("L86").Value = 5354.5
$ws.Range("M86").Value = -2499.2
$ws.Range("N86").Value = -7600.5

$ws.Range("H89").Value = 4392.1113
$ws.Range("I89").Value = 3622.2
$ws.Range("J89").Value = 5354.5
$ws.Range("K89").Value = 18111
$ws.Range("L89").Value = 26772.5
$ws.Range("M89").Value = -12495
$ws.Range("N89").Value = -38004.5

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").Value = $null

$ws.Range("H134").Value = 2315.2222
$ws.Range("I134").Value = 2305.2856
$ws.Range("K134").Value = 6915.8568
$ws.Range("M134").Value = -4380.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 27687.77
$ws.Range("J86").Value = 56126
$ws.Range("L86").Value = 56126
$ws.Range("N86").Value = -58372

$ws.Range("H89").Value = 27687.77
$ws.Range("J89").Value = 56126
$ws.Range("L89").Value = 280630
$ws.Range("N89").Value = -291862

$ws.Range("H132").Value = 2754.4443
$ws.Range("I132").Value = 1758.8
$ws.Range("K132").Value = 5276.4
$ws.Range("M132").Value = -2746.4

$ws.Range("H134").Value = 2294.25
$ws.Range("I134").Value = 2294.25
$ws.Range("K134").Value = 6882.75
$ws.Range("M134").Value = -4347.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 336260.66
$ws.Range("I80").Value = 4391.5
$ws.Range("K80").Value = 13174.5
$ws.Range("M80").Value = -12238.5

$ws.Range("H83").Value = 336260.66
$ws.Range("I83").Value = 4391.5
$ws.Range("K83").Value = 39523.5
$ws.Range("M83").Value = -34843.5

$ws.Range("H92").Value = 589.8
$ws.Range("I92").Value = 399.66666
$ws.Range("K92").Value = 1198.99998
$ws.Range("M92").Value = 49.00001999999995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3998.3333
$ws.Range("I80").Value = 2796.8
$ws.Range("K80").Value = 2796.8
$ws.Range("M80").Value = -1798.8

$ws.Range("H83").Value = 3998.3333
$ws.Range("I83").Value = 2796.8
$ws.Range("K83").Value = 13984
$ws.Range("M83").Value = -8992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 299.33334

$ws.Range("H61").Value = 4927.857
$ws.Range("I61").Value = 4900.2
$ws.Range("K61").Value = 4900.2
$ws.Range("M61").Value = -4698.2

$ws.Range("H93").Value = 2851.6
$ws.Range("I93").Value = 2842.75
$ws.Range("K93").Value = 2842.75
$ws.Range("M93").Value = -1594.75

$ws.Range("H113").Value = 4927.857
$ws.Range("I113").Value = 4900.2
$ws.Range("K113").Value = 4900.2
$ws.Range("M113").Value = -2730.2

$ws.Range("H132").Value = 3292.2273
$ws.Range("I132").Value = 2176.3635
$ws.Range("J132").Value = 4408.091
$ws.Range("K132").Value = 6529.0905
$ws.Range("L132").Value = 13224.273
$ws.Range("M132").Value = -3999.0905
$ws.Range("N132").Value = -18284.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2999
$ws.Range("I96").Value = 2999
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 2999
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -1626
$ws.Range("N96").Value = $null

$ws.Range("H107").Value = 350.375
$ws.Range("I107").Value = 369.57144
$ws.Range("J107").Value = 216
$ws.Range("K107").Value = 1108.71432
$ws.Range("L107").Value = 648
$ws.Range("M107").Value = 811.28568
$ws.Range("N107").Value = -4488

$ws.Range("H132").Value = 3519
$ws.Range("I132").Value = 2866.8572
$ws.Range("K132").Value = 8600.571599999999
$ws.Range("M132").Value = -6070.571599999999
